$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 65000
$ws.Range("J3").Value = 65000
$ws.Range("L3").Value = 65000
$ws.Range("N3").Value = -65228

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H57").Value = 54975
$ws.Range("I57").Value = 50000
$ws.Range("J57").Value = 59950
$ws.Range("K57").Value = 150000
$ws.Range("L57").Value = 179850
$ws.Range("M57").Value = -149501
$ws.Range("N57").Value = -180848

$ws.Range("H70").Value = 8845.714
$ws.Range("I70").Value = 14535.733
$ws.Range("J70").Value = 2280.3076
$ws.Range("K70").Value = 43607.199
$ws.Range("L70").Value = 6840.9228
$ws.Range("M70").Value = -43337.199
$ws.Range("N70").Value = -7380.9228

$ws.Range("H73").Value = 8845.714
$ws.Range("I73").Value = 14535.733
$ws.Range("J73").Value = 2280.3076
$ws.Range("K73").Value = 43607.199
$ws.Range("L73").Value = 6840.9228
$ws.Range("M73").Value = -42671.199
$ws.Range("N73").Value = -8712.9228

$ws.Range("H100").Value = 31173.314
$ws.Range("I100").Value = 35755.7
$ws.Range("K100").Value = 35755.7
$ws.Range("M100").Value = -35214.7

$ws.Range("H102").Value = 65000
$ws.Range("J102").Value = 65000
$ws.Range("L102").Value = 65000
$ws.Range("N102").Value = -71490

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2918.9443
$ws.Range("I32").Value = 2198.1077
$ws.Range("J32").Value = 9612.429
$ws.Range("K32").Value = 2198.1077
$ws.Range("L32").Value = 9612.429
$ws.Range("M32").Value = -1911.1077
$ws.Range("N32").Value = -10186.429

$ws.Range("H45").Value = 10667.6
$ws.Range("I45").Value = 13199.3
$ws.Range("J45").Value = 5604.2
$ws.Range("K45").Value = 13199.3
$ws.Range("L45").Value = 5604.2
$ws.Range("M45").Value = -12822.3
$ws.Range("N45").Value = -6358.2

$ws.Range("H60").Value = 50000
$ws.Range("J60").Value = 50000
$ws.Range("L60").Value = 50000
$ws.Range("N60").Value = -51466

$ws.Range("H102").Value = 4022.2104
$ws.Range("I102").Value = 4116.2856
$ws.Range("K102").Value = 4116.2856
$ws.Range("M102").Value = -2494.2856

$ws.Range("H122").Value = 1390.1428
$ws.Range("I122").Value = 1205.3334
$ws.Range("K122").Value = 3616.0002
$ws.Range("M122").Value = -1166.0002

$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 3832.625
$ws.Range("I132").Value = 3982.0527
$ws.Range("J132").Value = 3264.8
$ws.Range("K132").Value = 11946.1581
$ws.Range("L132").Value = 9794.400000000001
$ws.Range("M132").Value = -9416.158100000001
$ws.Range("N132").Value = -14854.4

$ws.Range("H137").Value = 75000
$ws.Range("J137").Value = 75000
$ws.Range("L137").Value = 75000
$ws.Range("N137").Value = -85200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 38749.5
$ws.Range("J58").Value = 37500
$ws.Range("L58").Value = 37500
$ws.Range("N58").Value = -38088

$ws.Range("H82").Value = 45116.145
$ws.Range("I82").Value = 30257
$ws.Range("J82").Value = 47592.668
$ws.Range("K82").Value = 30257
$ws.Range("L82").Value = 47592.668
$ws.Range("M82").Value = -29874
$ws.Range("N82").Value = -48358.668

$ws.Range("H85").Value = 45116.145
$ws.Range("I85").Value = 30257
$ws.Range("J85").Value = 47592.668
$ws.Range("K85").Value = 30257
$ws.Range("L85").Value = 47592.668
$ws.Range("M85").Value = -28931
$ws.Range("N85").Value = -50244.668

$ws.Range("H107").Value = 1587.0476
$ws.Range("I107").Value = 1566.45
$ws.Range("K107").Value = 1566.45
$ws.Range("M107").Value = 353.55

$ws.Range("H114").Value = 89999
$ws.Range("J114").Value = 89999
$ws.Range("L114").Value = 89999
$ws.Range("N114").Value = -98677

$ws.Range("H115").Value = 25000
$ws.Range("J115").Value = 25000
$ws.Range("L115").Value = 25000
$ws.Range("N115").Value = -28134

$ws.Range("H134").Value = 5881.386
$ws.Range("I134").Value = 3991.3333
$ws.Range("J134").Value = 11551.546
$ws.Range("K134").Value = 11973.9999
$ws.Range("L134").Value = 34654.638
$ws.Range("M134").Value = -9438.999899999999
$ws.Range("N134").Value = -39724.638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H35").Value = 360.33334
$ws.Range("I35").Value = 360.33334
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 360.33334
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -66.33334000000002
$ws.Range("N35").ClearContents()

$ws.Range("H58").Value = 3266.6924
$ws.Range("I58").Value = 1761.9
$ws.Range("K58").Value = 1761.9
$ws.Range("M58").Value = -1558.9

$ws.Range("H86").Value = 4999.4443
$ws.Range("J86").Value = 4998.75
$ws.Range("L86").Value = 4998.75
$ws.Range("N86").Value = -7244.75

$ws.Range("H89").Value = 4999.4443
$ws.Range("J89").Value = 4998.75
$ws.Range("L89").Value = 24993.75
$ws.Range("N89").Value = -36225.75

$ws.Range("H94").Value = 699.5833
$ws.Range("J94").Value = 842.2222
$ws.Range("L94").Value = 842.2222
$ws.Range("N94").Value = -1744.2222

$ws.Range("H105").Value = 2526.4285
$ws.Range("I105").Value = 2526.4285
$ws.Range("K105").Value = 2526.4285
$ws.Range("M105").Value = -779.4285

$ws.Range("H127").Value = 30000
$ws.Range("I127").Value = 30000
$ws.Range("K127").Value = 30000
$ws.Range("M127").Value = -25040

$ws.Range("H136").Value = 3266.6924
$ws.Range("I136").Value = 1761.9
$ws.Range("K136").Value = 5285.700000000001
$ws.Range("M136").Value = -2735.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3182.3333
$ws.Range("J34").Value = 5098.231
$ws.Range("L34").Value = 15294.693
$ws.Range("N34").Value = -15462.693

$ws.Range("H39").Value = 7249.7
$ws.Range("J39").Value = 7777.4443
$ws.Range("L39").Value = 23332.3329
$ws.Range("N39").Value = -23920.3329

$ws.Range("H55").Value = 1515.3529
$ws.Range("I55").Value = 572
$ws.Range("J55").Value = 1805.6154
$ws.Range("K55").Value = 1716
$ws.Range("L55").Value = 5416.8462
$ws.Range("M55").Value = -1539
$ws.Range("N55").Value = -5770.8462

$ws.Range("H61").Value = 1243.6
$ws.Range("I61").Value = 120
$ws.Range("J61").Value = 1992.6666
$ws.Range("K61").Value = 360
$ws.Range("L61").Value = 5977.9998
$ws.Range("M61").Value = -145
$ws.Range("N61").Value = -6407.9998

$ws.Range("H107").Value = 842.7
$ws.Range("J107").Value = 960.7778
$ws.Range("L107").Value = 2882.3334
$ws.Range("N107").Value = -6722.3334

$ws.Range("H131").Value = 16397.287
$ws.Range("I131").Value = 56339.168
$ws.Range("J131").Value = 1419.0834
$ws.Range("K131").Value = 169017.504
$ws.Range("L131").Value = 4257.2502
$ws.Range("M131").Value = -163977.504
$ws.Range("N131").Value = -14337.2502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 18667.334
$ws.Range("J26").Value = 18667.334
$ws.Range("L26").Value = 18667.334
$ws.Range("N26").Value = -19227.334

$ws.Range("H50").Value = 18667.334
$ws.Range("J50").Value = 18667.334
$ws.Range("L50").Value = 18667.334
$ws.Range("N50").Value = -19663.334

$ws.Range("H70").Value = 8333.75
$ws.Range("I70").Value = 7445.625
$ws.Range("J70").Value = 9221.875
$ws.Range("K70").Value = 7445.625
$ws.Range("L70").Value = 9221.875
$ws.Range("M70").Value = -7175.625
$ws.Range("N70").Value = -9761.875

$ws.Range("H73").Value = 8333.75
$ws.Range("I73").Value = 7445.625
$ws.Range("J73").Value = 9221.875
$ws.Range("K73").Value = 7445.625
$ws.Range("L73").Value = 9221.875
$ws.Range("M73").Value = -6509.625
$ws.Range("N73").Value = -11093.875

$ws.Range("H102").Value = 2889.8
$ws.Range("I102").Value = 2739.25
$ws.Range("K102").Value = 2739.25
$ws.Range("M102").Value = -1117.25

$ws.Range("H122").Value = 1469.1
$ws.Range("I122").Value = 1465.6666
$ws.Range("K122").Value = 4396.9998
$ws.Range("M122").Value = -1946.9998

$ws.Range("H132").Value = 2980.7878
$ws.Range("I132").Value = 3076.5925
$ws.Range("J132").Value = 2549.6667
$ws.Range("K132").Value = 9229.7775
$ws.Range("L132").Value = 7649.000100000001
$ws.Range("M132").Value = -6699.7775
$ws.Range("N132").Value = -12709.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 9299.666999999999
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9299.666999999999
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 9299.666999999999
$ws.Range("N34").Value = -9643.666999999999
$ws.Range("M34").ClearContents()

$ws.Range("H61").Value = 1811
$ws.Range("I61").Value = 1558.85
$ws.Range("K61").Value = 1558.85
$ws.Range("M61").Value = -1356.85

$ws.Range("H113").Value = 1811
$ws.Range("I113").Value = 1558.85
$ws.Range("K113").Value = 1558.85
$ws.Range("M113").Value = 611.1500000000001

$ws.Range("H132").Value = 2702
$ws.Range("I132").Value = 2465.348
$ws.Range("K132").Value = 7396.044
$ws.Range("M132").Value = -4866.044

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H132").Value = 30346.879
$ws.Range("I132").Value = 26504.77
$ws.Range("K132").Value = 79514.31
$ws.Range("M132").Value = -76984.31
